# diverse Updates, einige Bilder hinzugefügt
#
# The "Title" column (H) for several recurring TK events previously stored a
# redundant "TK - <Veranstaltung>" string; the prefix is dropped so the
# Title now simply mirrors the event name (re-using the already existing
# shared string). The "Musikfest Oberstdorf" entry instead gets a new,
# more descriptive title. A new event (row 62, "Erstkommunion") also gains
# the four extra metadata columns (Öffentlich / Beschreibung / Eintritt /
# Kategorie) that the other events in this block already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62 - Erstkommunion: drop the "TK - " prefix from the title and add
# the Öffentlich/Beschreibung/Eintritt/Kategorie metadata columns.
$ws.Range("H62").Value = "Erstkommunion"
$ws.Range("M62").Value = "Ja"
$ws.Range("N62").Value = "Testbeschreibung"
$ws.Range("O62").Value = "frei"
$ws.Range("P62").Value = "krichliches"

# Row 65 - Jahreskonzert
$ws.Range("H65").Value = "Jahreskonzert"

# Row 66 - Firmung
$ws.Range("H66").Value = "Firmung"

# Row 67 - 210 Jahre TK Riezlern
$ws.Range("H67").Value = "210 Jahre TK Riezlern"

# Rows 68-76 - Sommerkonzert (several dates) and Tag der Blasmusik
$ws.Range("H68").Value = "Sommerkonzert"
$ws.Range("H69").Value = "Sommerkonzert"
$ws.Range("H70").Value = "Sommerkonzert"
$ws.Range("H71").Value = "Tag der Blasmusik"
$ws.Range("H72").Value = "Sommerkonzert"

# Row 73 - Musikfest Oberstdorf gets a brand new title
$ws.Range("H73").Value = "TK Riezlern @ Musikfest Oberstdorf"

$ws.Range("H74").Value = "Sommerkonzert"
$ws.Range("H75").Value = "Sommerkonzert"
$ws.Range("H76").Value = "Sommerkonzert"

# Row 77 - Kriegergedenken
$ws.Range("H77").Value = "Kriegergedenken"

# Update the visible selection to reflect where the edits were made.
$ws.Activate()
$ws.Range("S76").Select()
